# Applies the "cryptos list" price/volume refresh captured in the commit
# "Updated cryptos list on Sat Nov 11 21:38:05 UTC 2023 with GitHub Actions".
# Only cell VALUES change (Coin/Link/Price/Volume(1h) columns); no rows are
# inserted or removed. Price cells that look like plain numbers are protected
# with a text NumberFormat ("@") before the write so Excel keeps them as the
# exact original text (e.g. "0.670", "252.31") instead of silently coercing
# them into numeric values and dropping formatting/trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "37.044.99"
$ws.Range("E2").Value = "  -0.93%  "
# Row 3
$ws.Range("D3").Value = "2.052.68"
$ws.Range("E3").Value = "  -2.20%  "
# Row 4
$ws.Range("E4").Value = "  +0.26%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.31"
$ws.Range("E5").Value = "  +0.03%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.670"
$ws.Range("E6").Value = "  +0.30%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.55"
$ws.Range("E7").Value = "  +7.66%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.22"
$ws.Range("E9").Value = "  -1.93%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.386"
$ws.Range("E10").Value = "  +1.34%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0796"
$ws.Range("E11").Value = "  +5.36%  "
# Row 12
$ws.Range("E12").Value = "  +1.96%  "
# Row 13
$ws.Range("E13").Value = "  +8.26%  "
# Row 14
$ws.Range("D14").Value = "2.355.03"
$ws.Range("E14").Value = "  -1.97%  "
# Row 15
$ws.Range("E15").Value = "  -3.76%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.61"
$ws.Range("E16").Value = "  +7.03%  "
# Row 17
$ws.Range("D17").Value = "2.064.10"
$ws.Range("E17").Value = "  -1.61%  "
# Row 18
$ws.Range("D18").Value = "36.993.51"
$ws.Range("E18").Value = "  -0.78%  "
# Row 19
$ws.Range("E19").Value = "  +14.16%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "75.53"
$ws.Range("E20").Value = "  +3.08%  "
# Row 21
$ws.Range("D21").Value = "0.0₃0912"
$ws.Range("E21").Value = "  +6.62%  "
# Row 22
$ws.Range("E22").Value = "  +3.85%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.00"
$ws.Range("E23").Value = "  -1.53%  "
# Row 24
$ws.Range("E24").Value = "  -0.02%  "
# Row 25
$ws.Range("E25").Value = "  -3.40%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("E26").Value = "  +11.50%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.26"
$ws.Range("E27").Value = "  -1.55%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.28"
$ws.Range("E28").Value = "  -0.11%  "
# Row 29
$ws.Range("E29").Value = "  -3.30%  "
# Row 30
$ws.Range("E30").Value = "  +1.47%  "
# Row 31
$ws.Range("E31").Value = "  +4.84%  "
# Row 32
$ws.Range("E32").Value = "  +4.51%  "
# Row 33
$ws.Range("E33").Value = "  -0.42%  "
# Row 34
$ws.Range("E34").Value = "  +5.90%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0884"
$ws.Range("E35").Value = "  -2.43%  "
# Row 36
$ws.Range("E36").Value = "  +0.10%  "
# Row 37
$ws.Range("E37").Value = "  -0.10%  "
# Row 38
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.76"
$ws.Range("E38").Value = "  -3.14%  "
# Row 39
$ws.Range("B39").Value = "Cronos"
$ws.Range("C39").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.115"
$ws.Range("E39").Value = "  +18.47%  "
# Row 40
$ws.Range("E40").Value = "  +1.08%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.95"
$ws.Range("E41").Value = "  -1.42%  "
# Row 42
$ws.Range("E42").Value = "  -2.00%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.15"
$ws.Range("E43").Value = "  -2.49%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.43"
$ws.Range("E44").Value = "  -2.10%  "
# Row 45
$ws.Range("E45").Value = "  +1.98%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.52"
$ws.Range("E46").Value = "  +10.86%  "
# Row 47
$ws.Range("E47").Value = "  +5.00%  "
# Row 48
$ws.Range("D48").Value = "1.291.93"
$ws.Range("E48").Value = "  -3.17%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.90"
$ws.Range("E49").Value = "  -1.54%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.89"
$ws.Range("E50").Value = "  -1.11%  "
# Row 51
$ws.Range("D51").Value = "2.245.38"
$ws.Range("E51").Value = "  -1.83%  "

Write-Host "Applied cryptos list update"
